$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 5 (pushing the existing rows 5..63 down to 6..64)
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with this week's new data point (latest date)
$ws.Cells.Item(5,1).Value = 11
$ws.Cells.Item(5,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(5,3).Value = "Bíobío"
$ws.Cells.Item(5,4).Value = 45035
$ws.Cells.Item(5,5).Value = 8
$ws.Cells.Item(5,6).Value = 100112030
$ws.Cells.Item(5,7).Value = "Poroto granado"
$ws.Cells.Item(5,8).Value = "Sin especificar"
$ws.Cells.Item(5,9).Value = "Primera"
$ws.Cells.Item(5,10).Value = 100
$ws.Cells.Item(5,11).Value = 33000
$ws.Cells.Item(5,12).Value = 35000
$ws.Cells.Item(5,13).Value = 34000
$ws.Cells.Item(5,14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(5,15).Value = "Región Metropolitana"
$ws.Cells.Item(5,16).Value = 1360
$ws.Cells.Item(5,17).Value = 25
$ws.Cells.Item(5,18).Value = "Hortaliza"
